$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("E2").Value = 23.78000000000028
$ws.Range("H2").Value = [double]"1.198621349122976e-16"
$ws.Range("K2").Value = 46.67452524640071
$ws.Range("L2").Value = "[44.278604055658086, 49.070446437143325]"
$ws.Range("O2").Value = 1.540921321580579
$ws.Range("P2").Value = "[1.490605523324887, 1.5912371198362711]"
$ws.Range("S2").Value = 52.32068771904981
$ws.Range("T2").Value = "[50.687754870801356, 53.95362056729827]"
$ws.Range("W2").Value = 17.94806806806828
$ws.Range("X2").Value = 17.75763763763785
$ws.Range("Y2").Value = 18.13849849849871

# --- Row 3 updates ---
$ws.Range("E3").Value = 23.71000000000027
$ws.Range("H3").Value = [double]"1.198621349122976e-16"
$ws.Range("K3").Value = 47.93473820313206
$ws.Range("L3").Value = "[43.32420684418271, 52.5452695620814]"
$ws.Range("O3").Value = 1.13839493553504
$ws.Range("P3").Value = "[1.0377633390236554, 1.2390265320464255]"
$ws.Range("S3").Value = 51.10796685800172
$ws.Range("T3").Value = "[48.59759522107552, 53.61833849492791]"
$ws.Range("W3").Value = 19.41419419419441
$ws.Range("X3").Value = 19.03445445445467
$ws.Range("Y3").Value = 19.79393393393416
